$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.730.80"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.849.72"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("D4").Value = "'1.014"
$ws.Range("E4").Value = "  -2.60%  "
$ws.Range("D5").Value = "'319.98"
$ws.Range("D6").Value = "'1.012"
$ws.Range("E6").Value = "  -2.38%  "
$ws.Range("D7").Value = "'0.4329"
$ws.Range("E7").Value = "  -2.19%  "
$ws.Range("D8").Value = "'0.3772"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("D9").Value = "'0.07391"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").Value = "'0.8837"
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("D11").Value = "'21.64"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "1.855.78"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").Value = "'6.759"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("D15").Value = "'0.07130"
$ws.Range("E15").Value = "  -1.58%  "
$ws.Range("D16").Value = "'88.39"
$ws.Range("E16").Value = "  +5.42%  "
$ws.Range("E17").Value = "  -2.49%  "
$ws.Range("D18").Value = "'0.000009040"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("E19").Value = "  -2.40%  "
$ws.Range("D20").Value = "'15.55"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "27.726.36"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").Value = "'5.269"
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("D23").Value = "'11.21"
$ws.Range("E23").Value = "  -1.46%  "
$ws.Range("D24").Value = "2.091.11"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'2.026"
$ws.Range("E25").Value = "  +2.73%  "
$ws.Range("D26").Value = "'155.77"
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("E27").Value = "  -1.42%  "
$ws.Range("D28").Value = "'2.134"
$ws.Range("E28").Value = "  +7.21%  "
$ws.Range("D29").Value = "'5.433"
$ws.Range("E29").Value = "  +2.01%  "
$ws.Range("D30").Value = "'120.65"
$ws.Range("E30").Value = "  +2.44%  "
$ws.Range("D31").Value = "'0.08958"
$ws.Range("E31").Value = "  -1.70%  "
$ws.Range("D32").Value = "'1.241"
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("D33").Value = "'0.7799"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").Value = "'4.580"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("D35").Value = "'2.921"
$ws.Range("E35").Value = "  -4.41%  "
$ws.Range("D36").Value = "'1.148"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").Value = "'1.012"
$ws.Range("E37").Value = "  -2.44%  "
$ws.Range("D38").Value = "'0.05339"
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("D39").Value = "'0.01973"
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("D40").Value = "'7.169"
$ws.Range("E40").Value = "  +3.78%  "
$ws.Range("D41").Value = "'2.869"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").Value = "'0.5194"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").Value = "'8.961"
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("D45").Value = "'110.95"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("E46").Value = "  +0.79%  "
$ws.Range("D47").Value = "'1.719"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("D48").Value = "'0.4749"
$ws.Range("E48").Value = "  +0.59%  "
$ws.Range("D49").Value = "'0.06513"
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("D50").Value = "'1.013"
$ws.Range("E50").Value = "  -2.54%  "
$ws.Range("D51").Value = "'1.902"
$ws.Range("E51").Value = "  +0.27%  "
